$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.216.75"
Set-TextValue $ws.Range("E2") "  +0.05%  "

Set-TextValue $ws.Range("D3") "1.860.23"
Set-TextValue $ws.Range("E3") "  -0.07%  "

Set-TextValue $ws.Range("E4") "  -0.05%  "

Set-TextValue $ws.Range("D5") "235.97"
Set-TextValue $ws.Range("E5") "  +0.58%  "

Set-TextValue $ws.Range("E6") "  -0.11%  "

Set-TextValue $ws.Range("E7") "  +1.31%  "

Set-TextValue $ws.Range("E8") "  +2.69%  "

Set-TextValue $ws.Range("D9") "0.06562"
Set-TextValue $ws.Range("E9") "  +0.61%  "

Set-TextValue $ws.Range("D10") "21.80"
Set-TextValue $ws.Range("E10") "  +3.01%  "

Set-TextValue $ws.Range("D11") "0.07938"
Set-TextValue $ws.Range("E11") "  +1.15%  "

Set-TextValue $ws.Range("D12") "97.77"
Set-TextValue $ws.Range("E12") "  +0.74%  "

Set-TextValue $ws.Range("D13") "1.868.77"
Set-TextValue $ws.Range("E13") "  +0.29%  "

Set-TextValue $ws.Range("D14") "5.132"
Set-TextValue $ws.Range("E14") "  +0.91%  "

Set-TextValue $ws.Range("D15") "0.6804"
Set-TextValue $ws.Range("E15") "  +1.55%  "

Set-TextValue $ws.Range("D16") "266.90"
Set-TextValue $ws.Range("E16") "  -4.29%  "

Set-TextValue $ws.Range("D17") "30.216.22"
Set-TextValue $ws.Range("E17") "  +0.08%  "

Set-TextValue $ws.Range("D18") "13.70"
Set-TextValue $ws.Range("E18") "  +8.60%  "

Set-TextValue $ws.Range("D19") "0.000007575"
Set-TextValue $ws.Range("E19") "  +4.19%  "

Set-TextValue $ws.Range("D20") "1.001"

Set-TextValue $ws.Range("D21") "2.107.89"
Set-TextValue $ws.Range("E21") "  -0.34%  "

Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  -0.08%  "

Set-TextValue $ws.Range("D23") "5.259"
Set-TextValue $ws.Range("E23") "  -4.32%  "

Set-TextValue $ws.Range("D24") "6.173"
Set-TextValue $ws.Range("E24") "  +0.60%  "

Set-TextValue $ws.Range("D25") "167.25"
Set-TextValue $ws.Range("E25") "  +1.52%  "

Set-TextValue $ws.Range("D26") "9.189"
Set-TextValue $ws.Range("E26") "  +0.13%  "

Set-TextValue $ws.Range("D27") "18.93"
Set-TextValue $ws.Range("E27") "  -0.70%  "

Set-TextValue $ws.Range("D28") "1.949"
Set-TextValue $ws.Range("E28") "  +1.70%  "

Set-TextValue $ws.Range("D29") "1.396"
Set-TextValue $ws.Range("E29") "  +0.94%  "

Set-TextValue $ws.Range("D30") "0.09884"
Set-TextValue $ws.Range("E30") "  +2.18%  "

Set-TextValue $ws.Range("D31") "4.336"
Set-TextValue $ws.Range("E31") "  -1.60%  "

Set-TextValue $ws.Range("D32") "1.468"
Set-TextValue $ws.Range("E32") "  -0.34%  "

Set-TextValue $ws.Range("D33") "4.012"
Set-TextValue $ws.Range("E33") "  -1.39%  "

Set-TextValue $ws.Range("D34") "0.04707"
Set-TextValue $ws.Range("E34") "  +0.70%  "

Set-TextValue $ws.Range("D35") "1.128"
Set-TextValue $ws.Range("E35") "  +1.53%  "

Set-TextValue $ws.Range("D36") "0.7011"
Set-TextValue $ws.Range("E36") "  -0.12%  "

Set-TextValue $ws.Range("D37") "2.706"
Set-TextValue $ws.Range("E37") "  -0.86%  "

Set-TextValue $ws.Range("E38") "  +1.94%  "

Set-TextValue $ws.Range("D39") "2.618"
Set-TextValue $ws.Range("E39") "  +3.45%  "

Set-TextValue $ws.Range("D40") "6.331"
Set-TextValue $ws.Range("E40") "  +1.42%  "

Set-TextValue $ws.Range("D41") "73.82"
Set-TextValue $ws.Range("E41") "  +0.84%  "

Set-TextValue $ws.Range("E42") "  +0.10%  "

Set-TextValue $ws.Range("D43") "0.8412"
Set-TextValue $ws.Range("E43") "  -0.44%  "

Set-TextValue $ws.Range("E44") "  -0.14%  "

Set-TextValue $ws.Range("D45") "0.4157"
Set-TextValue $ws.Range("E45") "  +0.21%  "

Set-TextValue $ws.Range("D46") "103.15"
Set-TextValue $ws.Range("E46") "  -0.67%  "

Set-TextValue $ws.Range("E47") "  -0.16%  "

Set-TextValue $ws.Range("D48") "942.98"
Set-TextValue $ws.Range("E48") "  +0.77%  "

Set-TextValue $ws.Range("D49") "9.163"
Set-TextValue $ws.Range("E49") "  +0.75%  "

Set-TextValue $ws.Range("D50") "34.11"
Set-TextValue $ws.Range("E50") "  +0.54%  "

Set-TextValue $ws.Range("E51") "  +0.66%  "

